# Weekly fruit/vegetable price update: shift existing rows down by one week
# (each row absorbs the prior row's data) and prepend a new latest-week
# entry in row 13, pushing the oldest entry into a new row 38.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("D13").Value = 44526
$ws.Range("K13").Value = 1700
$ws.Range("L13").Value = 1700

# Row 14
$ws.Range("D14").Value = 44519
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 1600
$ws.Range("M14").Value = 1700
$ws.Range("P14").Value = 1700

# Row 15
$ws.Range("D15").Value = 44497
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 1800
$ws.Range("L15").Value = 1800
$ws.Range("M15").Value = 1800
$ws.Range("P15").Value = 1800

# Row 16
$ws.Range("D16").Value = 44482
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 2000
$ws.Range("P16").Value = 2000

# Row 17
$ws.Range("D17").Value = 44516
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("J17").Value = 360

# Row 18
$ws.Range("D18").Value = 44168
$ws.Range("H18").Value = 'Verde'
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1600
$ws.Range("L18").Value = 1600
$ws.Range("M18").Value = 1600
$ws.Range("P18").Value = 1600

# Row 19
$ws.Range("D19").Value = 44475
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 2000
$ws.Range("P19").Value = 2000

# Row 20
$ws.Range("D20").Value = 44162
$ws.Range("H20").Value = 'Verde'
$ws.Range("J20").Value = 700
$ws.Range("K20").Value = 1600
$ws.Range("L20").Value = 1600
$ws.Range("M20").Value = 1600
$ws.Range("P20").Value = 1600

# Row 21
$ws.Range("D21").Value = 44496
$ws.Range("J21").Value = 84
$ws.Range("K21").Value = 1800
$ws.Range("L21").Value = 1800
$ws.Range("M21").Value = 1800
$ws.Range("P21").Value = 1800

# Row 22
$ws.Range("D22").Value = 44498
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1600
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = 1600
$ws.Range("P22").Value = 1600

# Row 23
$ws.Range("D23").Value = 44487
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 1800
$ws.Range("L23").Value = 1800
$ws.Range("M23").Value = 1800
$ws.Range("P23").Value = 1800

# Row 24
$ws.Range("D24").Value = 44509
$ws.Range("J24").Value = 550

# Row 25
$ws.Range("D25").Value = 44494
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 1700
$ws.Range("L25").Value = 1700
$ws.Range("M25").Value = 1700
$ws.Range("P25").Value = 1700

# Row 26
$ws.Range("D26").Value = 44518
$ws.Range("J26").Value = 180
$ws.Range("K26").Value = 1600
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = 1600
$ws.Range("P26").Value = 1600

# Row 27
$ws.Range("D27").Value = 44481
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = 1850
$ws.Range("P27").Value = 1850

# Row 28
$ws.Range("D28").Value = 44523
$ws.Range("H28").Value = 'Sin especificar'
$ws.Range("J28").Value = 520
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = 1800
$ws.Range("P28").Value = 1800

# Row 29
$ws.Range("D29").Value = 44169
$ws.Range("H29").Value = 'Verde'
$ws.Range("J29").Value = 600
$ws.Range("K29").Value = 1600
$ws.Range("L29").Value = 1600
$ws.Range("M29").Value = 1600
$ws.Range("P29").Value = 1600

# Row 30
$ws.Range("D30").Value = 44474
$ws.Range("H30").Value = 'Sin especificar'
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 2000
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = 2000
$ws.Range("P30").Value = 2000

# Row 31
$ws.Range("D31").Value = 44176
$ws.Range("H31").Value = 'Verde'
$ws.Range("J31").Value = 700

# Row 32
$ws.Range("D32").Value = 44504
$ws.Range("K32").Value = 1600
$ws.Range("L32").Value = 1600
$ws.Range("M32").Value = 1600
$ws.Range("P32").Value = 1600

# Row 33
$ws.Range("D33").Value = 44522
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 1800
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = 1800
$ws.Range("P33").Value = 1800

# Row 34
$ws.Range("D34").Value = 44491
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("K34").Value = 1700
$ws.Range("L34").Value = 1700
$ws.Range("M34").Value = 1700
$ws.Range("P34").Value = 1700

# Row 35
$ws.Range("D35").Value = 44166
$ws.Range("H35").Value = 'Verde'
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 1600
$ws.Range("L35").Value = 1600
$ws.Range("M35").Value = 1600
$ws.Range("P35").Value = 1600

# Row 36
$ws.Range("D36").Value = 44495
$ws.Range("H36").Value = 'Sin especificar'
$ws.Range("J36").Value = 520
$ws.Range("K36").Value = 1800
$ws.Range("L36").Value = 1800
$ws.Range("M36").Value = 1800
$ws.Range("P36").Value = 1800

# Row 37
$ws.Range("D37").Value = 44161
$ws.Range("H37").Value = 'Verde'
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 1700
$ws.Range("L37").Value = 1700
$ws.Range("M37").Value = 1700
$ws.Range("P37").Value = 1700

# Row 38 (new row appended at the end)
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C38").Value = 'Los Lagos'
$ws.Range("D38").Value = 44517
$ws.Range("D38").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 300000000
$ws.Range("G38").Value = 'Espárragos'
$ws.Range("H38").Value = 'Sin especificar'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 120
$ws.Range("K38").Value = 1600
$ws.Range("L38").Value = 1600
$ws.Range("M38").Value = 1600
$ws.Range("N38").Value = '$/kilo'
$ws.Range("O38").Value = 'Provincia de Linares'
$ws.Range("P38").Value = 1600
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = 'Hortaliza'

